# Auto-generated edit script: updates Leve market-price / profit figures
# across multiple job sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR) per the
# scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1473
$ws.Range("I18").Value = 366.25
$ws.Range("K18").Value = 366.25
$ws.Range("M18").Value = -82.25
$ws.Range("H41").Value = 96.69231000000001
$ws.Range("I41").Value = 72.90909000000001
$ws.Range("J41").Value = 227.5
$ws.Range("K41").Value = 72.90909000000001
$ws.Range("L41").Value = 227.5
$ws.Range("M41").Value = 367.09091
$ws.Range("N41").Value = -1107.5
$ws.Range("H141").Value = 2431.25
$ws.Range("I141").Value = 1760
$ws.Range("J141").Value = 3102.5
$ws.Range("K141").Value = 5280
$ws.Range("L141").Value = 9307.5
$ws.Range("M141").Value = -100
$ws.Range("N141").Value = -19667.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7456.273
$ws.Range("I32").Value = 7062.48
$ws.Range("K32").Value = 7062.48
$ws.Range("M32").Value = -6775.48
$ws.Range("H45").Value = 6104.5835
$ws.Range("I45").Value = 5972.8887
$ws.Range("J45").Value = 6499.6665
$ws.Range("K45").Value = 5972.8887
$ws.Range("L45").Value = 6499.6665
$ws.Range("M45").Value = -5595.8887
$ws.Range("N45").Value = -7253.6665
$ws.Range("H108").Value = 20000
$ws.Range("I108").Value = 20000
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 20000
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -16160
$ws.Range("N108").ClearContents()
$ws.Range("H115").Value = 35000
$ws.Range("J115").Value = 35000
$ws.Range("L115").Value = 35000
$ws.Range("N115").Value = -38134
$ws.Range("H127").Value = 49980
$ws.Range("J127").Value = 49980
$ws.Range("L127").Value = 49980
$ws.Range("N127").Value = -59900
$ws.Range("H137").Value = 56987.5
$ws.Range("J137").Value = 56987.5
$ws.Range("L137").Value = 56987.5
$ws.Range("N137").Value = -67187.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 68520
$ws.Range("J50").Value = 68520
$ws.Range("L50").Value = 68520
$ws.Range("N50").Value = -69668
$ws.Range("H74").Value = 50000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51872
$ws.Range("H77").Value = 50000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159360
$ws.Range("H80").Value = 1921.1111
$ws.Range("J80").Value = 995.6667
$ws.Range("L80").Value = 995.6667
$ws.Range("N80").Value = -2991.6667
$ws.Range("H83").Value = 1921.1111
$ws.Range("J83").Value = 995.6667
$ws.Range("L83").Value = 4978.3335
$ws.Range("N83").Value = -14962.3335
$ws.Range("H127").Value = 59779.5
$ws.Range("J127").Value = 59779.5
$ws.Range("L127").Value = 59779.5
$ws.Range("N127").Value = -69699.5
$ws.Range("H130").Value = 65292.668
$ws.Range("J130").Value = 65292.668
$ws.Range("L130").Value = 65292.668
$ws.Range("N130").Value = -75332.66800000001
$ws.Range("H135").Value = 77755.39999999999
$ws.Range("J135").Value = 77755.39999999999
$ws.Range("L135").Value = 77755.39999999999
$ws.Range("N135").Value = -87895.39999999999
$ws.Range("H138").Value = 71926.5
$ws.Range("I138").Value = 71926.5
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 71926.5
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -66786.5
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 155616.17
$ws.Range("J53").Value = 175139.4
$ws.Range("L53").Value = 175139.4
$ws.Range("N53").Value = -176353.4
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -54492
$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -55242
$ws.Range("H117").Value = 12430
$ws.Range("J117").Value = 12000
$ws.Range("L117").Value = 12000
$ws.Range("N117").Value = -21178
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H121").Value = 55500
$ws.Range("I121").Value = 74000
$ws.Range("J121").Value = 37000
$ws.Range("K121").Value = 74000
$ws.Range("L121").Value = 37000
$ws.Range("M121").Value = -72690
$ws.Range("N121").Value = -39620
$ws.Range("H132").Value = 98696.95
$ws.Range("I132").Value = 98696.95
$ws.Range("K132").Value = 296090.85
$ws.Range("M132").Value = -293560.85
$ws.Range("H133").Value = 63756.145
$ws.Range("I133").Value = 31000
$ws.Range("J133").Value = 66275.84
$ws.Range("K133").Value = 31000
$ws.Range("L133").Value = 66275.84
$ws.Range("M133").Value = -28470
$ws.Range("N133").Value = -71335.84
$ws.Range("H134").Value = 2447.2
$ws.Range("I134").Value = 1659.1
$ws.Range("K134").Value = 4977.299999999999
$ws.Range("M134").Value = -2442.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2635.7144
$ws.Range("I5").Value = 362.5
$ws.Range("K5").Value = 1087.5
$ws.Range("M5").Value = -975.5
$ws.Range("H135").Value = 2635.7144
$ws.Range("I135").Value = 362.5
$ws.Range("K135").Value = 3262.5
$ws.Range("M135").Value = -727.5
$ws.Range("H138").Value = 1998.4286
$ws.Range("I138").Value = 1664.8334
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 4994.5002
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = 145.4997999999996
$ws.Range("N138").Value = -22280
$ws.Range("H141").Value = 30677.334
$ws.Range("I141").Value = 43000
$ws.Range("K141").Value = 129000
$ws.Range("M141").Value = -123820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4722.7
$ws.Range("I100").Value = 3604
$ws.Range("J100").Value = 7333
$ws.Range("K100").Value = 3604
$ws.Range("L100").Value = 7333
$ws.Range("M100").Value = -3063
$ws.Range("N100").Value = -8415
$ws.Range("H104").Value = 19123
$ws.Range("J104").Value = 19123
$ws.Range("L104").Value = 19123
$ws.Range("N104").Value = -26111
$ws.Range("H131").Value = 89977.336
$ws.Range("J131").Value = 89977.336
$ws.Range("L131").Value = 89977.336
$ws.Range("N131").Value = -100057.336
$ws.Range("H136").Value = 6443.4165
$ws.Range("I136").Value = 2887.5
$ws.Range("K136").Value = 8662.5
$ws.Range("M136").Value = -6112.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 22996.4
$ws.Range("J70").Value = 36500
$ws.Range("L70").Value = 36500
$ws.Range("N70").Value = -37130
$ws.Range("H73").Value = 22996.4
$ws.Range("J73").Value = 36500
$ws.Range("L73").Value = 36500
$ws.Range("N73").Value = -38684
$ws.Range("H75").Value = 42494.5
$ws.Range("J75").Value = 44989
$ws.Range("L75").Value = 44989
$ws.Range("N75").Value = -46861
$ws.Range("H78").Value = 42494.5
$ws.Range("J78").Value = 44989
$ws.Range("L78").Value = 134967
$ws.Range("N78").Value = -144327
$ws.Range("H86").Value = 12535612
$ws.Range("J86").Value = 32999.2
$ws.Range("L86").Value = 32999.2
$ws.Range("N86").Value = -35245.2
$ws.Range("H89").Value = 12535612
$ws.Range("J89").Value = 32999.2
$ws.Range("L89").Value = 164996
$ws.Range("N89").Value = -176228
$ws.Range("H93").Value = 71142.86
$ws.Range("I93").Value = 74000
$ws.Range("J93").Value = 70000
$ws.Range("K93").Value = 74000
$ws.Range("L93").Value = 70000
$ws.Range("M93").Value = -71504
$ws.Range("N93").Value = -74992
$ws.Range("H106").Value = 49975
$ws.Range("I106").Value = 49975
$ws.Range("K106").Value = 49975
$ws.Range("M106").Value = -48713
$ws.Range("H109").Value = 85833.164
$ws.Range("I109").Value = 76666.336
$ws.Range("J109").Value = 95000
$ws.Range("K109").Value = 76666.336
$ws.Range("L109").Value = 95000
$ws.Range("M109").Value = -75279.336
$ws.Range("N109").Value = -97774
$ws.Range("H125").Value = 61215.43
$ws.Range("J125").Value = 61215.43
$ws.Range("L125").Value = 61215.43
$ws.Range("N125").Value = -71055.42999999999
$ws.Range("H126").Value = 16674960
$ws.Range("I126").Value = 16674960
$ws.Range("K126").Value = 50024880
$ws.Range("M126").Value = -50022410
$ws.Range("H129").Value = 30428
$ws.Range("J129").Value = 30428
$ws.Range("L129").Value = 30428
$ws.Range("N129").Value = -40428
$ws.Range("H139").Value = 119600
$ws.Range("I139").Value = 99333.336
$ws.Range("J139").Value = 150000
$ws.Range("K139").Value = 99333.336
$ws.Range("L139").Value = 150000
$ws.Range("M139").Value = -94193.336
$ws.Range("N139").Value = -160280
